$wb = $excel.ActiveWorkbook

# Insert the new "AUTO_SANTEI_MST" worksheet right after "TEN_MST",
# which pushes PT_SANTEI_CONF / USER_MST / TEKIOU_BYOMEI_MST / BYOMEI_MST down.
$afterSheet = $wb.Worksheets.Item("TEN_MST")
$newSheet = $wb.Worksheets.Add($null, $afterSheet)
$newSheet.Name = "AUTO_SANTEI_MST"

# Header row
$newSheet.Range("A1").Value = "hp_id"
$newSheet.Range("B1").Value = "item_cd"
$newSheet.Range("C1").Value = "seq_no"
$newSheet.Range("D1").Value = "start_date"
$newSheet.Range("E1").Value = "end_date"
$newSheet.Range("F1").Value = "create_date"
$newSheet.Range("G1").Value = "create_id"
$newSheet.Range("H1").Value = "create_machine"
$newSheet.Range("I1").Value = "update_date"
$newSheet.Range("J1").Value = "update_id"
$newSheet.Range("K1").Value = "update_machine"

# Data row
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = 111015970
$newSheet.Range("C2").Value = 0
$newSheet.Range("D2").Value = 0
$newSheet.Range("E2").Value = 99999999
$newSheet.Range("G2").Value = 0
$newSheet.Range("I2").NumberFormat = "m/d/yy h:mm"
$newSheet.Range("J2").Value = 2

# Column widths (best-fit like) similar to the reference sheet
$newSheet.Columns.Item(2).ColumnWidth = 9.17
$newSheet.Columns.Item(9).ColumnWidth = 13.83

# Match the reference selection/view state of the new sheet
$newSheet.Range("F4").Select() | Out-Null

Write-Host "Inserted AUTO_SANTEI_MST sheet with data"
